# Updated run for publication: refresh per-position nucleotide frequency table (B2:X5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Cells.Item(2, 2).Value = 0.00609756097560976  # B2
$ws.Cells.Item(2, 3).Value = 0.981707317073171  # C2
$ws.Cells.Item(2, 4).Value = 0.00203252032520325  # D2
$ws.Cells.Item(2, 5).Value = 0.00203252032520325  # E2
$ws.Cells.Item(2, 6).Value = 0.016260162601626  # F2
$ws.Cells.Item(2, 7).Value = 0.00203252032520325  # G2
$ws.Cells.Item(2, 8).Value = 0.0121951219512195  # H2
$ws.Cells.Item(2, 9).Value = 0.973577235772358  # I2
$ws.Cells.Item(2, 10).Value = 0.024390243902439  # J2
$ws.Cells.Item(2, 11).Value = 0.951219512195122  # K2
$ws.Cells.Item(2, 12).Value = 0.961382113821138  # L2
$ws.Cells.Item(2, 14).Value = 0.0040650406504065  # N2
$ws.Cells.Item(2, 15).Value = 0  # O2
$ws.Cells.Item(2, 17).Value = 0.995934959349594  # Q2
$ws.Cells.Item(2, 18).Value = 0.0040650406504065  # R2
$ws.Cells.Item(2, 19).Value = 0.890243902439024  # S2
$ws.Cells.Item(2, 20).Value = 0.0609756097560976  # T2
$ws.Cells.Item(2, 21).Value = 0.0691056910569106  # U2
$ws.Cells.Item(2, 22).Value = 0.989837398373984  # V2
$ws.Cells.Item(2, 23).Value = 0.00609756097560976  # W2
$ws.Cells.Item(2, 24).Value = 0.00203252032520325  # X2

# row 3
$ws.Cells.Item(3, 2).Value = 0.0101626016260163  # B3
$ws.Cells.Item(3, 3).Value = 0.0040650406504065  # C3
$ws.Cells.Item(3, 4).Value = 0.0223577235772358  # D3
$ws.Cells.Item(3, 5).Value = 0.00609756097560976  # E3
$ws.Cells.Item(3, 6).Value = 0.00203252032520325  # F3
$ws.Cells.Item(3, 8).Value = 0.977642276422764  # H3
$ws.Cells.Item(3, 9).Value = 0.0101626016260163  # I3
$ws.Cells.Item(3, 10).Value = 0.0142276422764228  # J3
$ws.Cells.Item(3, 11).Value = 0.00203252032520325  # K3
$ws.Cells.Item(3, 12).Value = 0.0040650406504065  # L3
$ws.Cells.Item(3, 13).Value = 0.943089430894309  # M3
$ws.Cells.Item(3, 14).Value = 0.00609756097560976  # N3
$ws.Cells.Item(3, 16).Value = 0.991869918699187  # P3
$ws.Cells.Item(3, 17).Value = 0.0040650406504065  # Q3
$ws.Cells.Item(3, 18).Value = 0.99390243902439  # R3
$ws.Cells.Item(3, 19).Value = 0.0975609756097561  # S3
$ws.Cells.Item(3, 20).Value = 0.930894308943089  # T3
$ws.Cells.Item(3, 21).Value = 0.92479674796748  # U3
$ws.Cells.Item(3, 22).Value = 0.00203252032520325  # V3
$ws.Cells.Item(3, 23).Value = 0.0040650406504065  # W3
$ws.Cells.Item(3, 24).Value = 0.0040650406504065  # X3

# row 4
$ws.Cells.Item(4, 2).Value = 0.979674796747967  # B4
$ws.Cells.Item(4, 3).Value = 0.00203252032520325  # C4
$ws.Cells.Item(4, 4).Value = 0.00813008130081301  # D4
$ws.Cells.Item(4, 5).Value = 0.0040650406504065  # E4
$ws.Cells.Item(4, 6).Value = 0.975609756097561  # F4
$ws.Cells.Item(4, 7).Value = 0.99390243902439  # G4
$ws.Cells.Item(4, 8).Value = 0.00813008130081301  # H4
$ws.Cells.Item(4, 9).Value = 0.0040650406504065  # I4
$ws.Cells.Item(4, 10).Value = 0.951219512195122  # J4
$ws.Cells.Item(4, 11).Value = 0.0040650406504065  # K4
$ws.Cells.Item(4, 12).Value = 0.032520325203252  # L4
$ws.Cells.Item(4, 13).Value = 0.00203252032520325  # M4
$ws.Cells.Item(4, 16).Value = 0.00203252032520325  # P4
$ws.Cells.Item(4, 18).Value = 0.00203252032520325  # R4
$ws.Cells.Item(4, 19).Value = 0.0101626016260163  # S4
$ws.Cells.Item(4, 20).Value = 0.00203252032520325  # T4
$ws.Cells.Item(4, 22).Value = 0.00609756097560976  # V4
$ws.Cells.Item(4, 23).Value = 0.98780487804878  # W4
$ws.Cells.Item(4, 24).Value = 0.991869918699187  # X4

# row 5
$ws.Cells.Item(5, 2).Value = 0.0040650406504065  # B5
$ws.Cells.Item(5, 3).Value = 0.0121951219512195  # C5
$ws.Cells.Item(5, 4).Value = 0.967479674796748  # D5
$ws.Cells.Item(5, 5).Value = 0.98780487804878  # E5
$ws.Cells.Item(5, 6).Value = 0.00609756097560976  # F5
$ws.Cells.Item(5, 7).Value = 0.0040650406504065  # G5
$ws.Cells.Item(5, 8).Value = 0.00203252032520325  # H5
$ws.Cells.Item(5, 10).Value = 0.0101626016260163  # J5
$ws.Cells.Item(5, 11).Value = 0.040650406504065  # K5
$ws.Cells.Item(5, 13).Value = 0.0528455284552846  # M5
$ws.Cells.Item(5, 14).Value = 0.989837398373984  # N5
$ws.Cells.Item(5, 16).Value = 0.00609756097560976  # P5
$ws.Cells.Item(5, 19).Value = 0.00203252032520325  # S5
$ws.Cells.Item(5, 20).Value = 0.0040650406504065  # T5
$ws.Cells.Item(5, 21).Value = 0.00609756097560976  # U5
$ws.Cells.Item(5, 22).Value = 0.00203252032520325  # V5
$ws.Cells.Item(5, 23).Value = 0.00203252032520325  # W5
$ws.Cells.Item(5, 24).Value = 0.00203252032520325  # X5
